$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to keep text-formatted numeric-looking values
# as literal text (matching the source data which stores these as strings),
# instead of letting Excel auto-convert them to numbers.

$ws.Range('D2').Value = '27.342.82'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '1.711.48'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.69'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.06691'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2670'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.92'
$ws.Range('E10').Value = '  -3.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07677'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.516'
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('D13').Value = '1.948.94'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').Value = '1.716.51'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5847'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').Value = '0.0₅8252'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '27.400.84'
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '222.50'
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.643'
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.42'
$ws.Range('E22').Value = '  -2.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.016'
$ws.Range('E23').Value = '  -1.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.68'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.688'
$ws.Range('E26').Value = '  -3.02%  '
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.256'
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.23'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('E30').Value = '  -4.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.294'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.460'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.439'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.641'
$ws.Range('E34').Value = '  -1.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.864'
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9528'
$ws.Range('E36').Value = '  -0.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.393'
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5872'
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01640'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('D40').Value = '1.090.77'
$ws.Range('E40').Value = '  +3.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.816'
$ws.Range('E41').Value = '  -1.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.004'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8436'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.05'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').Value = '1.855.02'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('D46').Value = '0.0₈112'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.93'
$ws.Range('E47').Value = '  -1.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4526'
$ws.Range('E48').Value = '  +2.07%  '
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.138'
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05234'
$ws.Range('E51').Value = '  -0.36%  '
